# Update the "Förändrad" date column (C) for rows 2-27 from 2023-12-30 (45290)
# to 2023-12-31 (45291), matching the automatic update of the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45290) {
        $cell.Value = 45291
    }
}
